$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 1
$ws.Range("G7").Value = 1
$ws.Range("G8").Value = 0
$ws.Range("G9").Value = 1
$ws.Range("G10").Value = 0
$ws.Range("G11").Value = 1
$ws.Range("G12").Value = 2
$ws.Range("G13").Value = 0
$ws.Range("G14").Value = 1
$ws.Range("G16").Value = 2
$ws.Range("G18").Value = 1
$ws.Range("G19").Value = 2
